# Natmi following Dr Hou advice
# Adds an "ECs" sending/target cluster to the Efna5-Epha3 LR-pair table,
# expanding the 2x2 FAPs/sCs cluster grid into a 3x3 ECs/FAPs/sCs grid
# with refreshed NATMI statistics for every cluster pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @{A="ECs"; B="Efna5"; C="Epha3"; D="ECs"; E=2; F=0.6666666666666666; G=0.4223226666666666; H=1.266968; I=0.1676547342089268; J=0.1676547342089269; K=2; L=0.6666666666666666; M=1.106217; N=3.318651; O=0.01813136626967656; P=0.01813136626967656; Q=0.4671805133519999; R=4.204624620168; S=0.003039809392787325; T=0.003039809392787326},
    @{A="ECs"; B="Efna5"; C="Epha3"; D="FAPs"; E=2; F=0.6666666666666666; G=0.4223226666666666; H=1.266968; I=0.1676547342089268; J=0.1676547342089269; K=3; L=1; M=59.24481466666666; N=177.734444; O=0.971047664518299; P=0.971047664518299; Q=25.02042811619911; R=225.183853045792; S=0.1628007380990146; T=0.1628007380990146},
    @{A="ECs"; B="Efna5"; C="Epha3"; D="sCs"; E=2; F=0.6666666666666666; G=0.4223226666666666; H=1.266968; I=0.1676547342089268; J=0.1676547342089269; K=3; L=1; M=0.6602006666666667; N=1.980602; O=0.01082096921202439; P=0.01082096921202439; Q=0.2788177060817777; R=2.509359354736; S=0.001814186717124929; T=0.001814186717124929},
    @{A="FAPs"; B="Efna5"; C="Epha3"; D="ECs"; E=3; F=1; G=1.874986333333333; H=5.624959; I=0.7443368783435028; J=0.7443368783435029; K=2; L=0.6666666666666666; M=1.106217; N=3.318651; O=0.01813136626967656; P=0.01813136626967656; Q=2.074141756701; R=18.667275810309; S=0.01349584456927373; T=0.01349584456927374},
    @{A="FAPs"; B="Efna5"; C="Epha3"; D="FAPs"; E=3; F=1; G=1.874986333333333; H=5.624959; I=0.7443368783435028; J=0.7443368783435029; K=3; L=1; M=59.24481466666666; N=177.734444; O=0.971047664518299; P=0.971047664518299; Q=111.0832178208662; R=999.748960387796; S=0.7227865873302997; T=0.7227865873302998},
    @{A="FAPs"; B="Efna5"; C="Epha3"; D="sCs"; E=3; F=1; G=1.874986333333333; H=5.624959; I=0.7443368783435028; J=0.7443368783435029; K=3; L=1; M=0.6602006666666667; N=1.980602; O=0.01082096921202439; P=0.01082096921202439; Q=1.237867227257555; R=11.140805045318; S=0.008054446443929384; T=0.008054446443929386},
    @{A="sCs"; B="Efna5"; C="Epha3"; D="ECs"; E=3; F=1; G=0.2216933333333333; H=0.66508; I=0.08800838744757017; J=0.08800838744757018; K=2; L=0.6666666666666666; M=1.106217; N=3.318651; O=0.01813136626967656; P=0.01813136626967656; Q=0.24524093412; R=2.20716840708; S=0.0015957123076155; T=0.0015957123076155},
    @{A="sCs"; B="Efna5"; C="Epha3"; D="FAPs"; E=3; F=1; G=0.2216933333333333; H=0.66508; I=0.08800838744757017; J=0.08800838744757018; K=3; L=1; M=59.24481466666666; N=177.734444; O=0.971047664518299; P=0.971047664518299; Q=13.13418044616889; R=118.20762401552; S=0.0854603390889846; T=0.08546033908898461},
    @{A="sCs"; B="Efna5"; C="Epha3"; D="sCs"; E=3; F=1; G=0.2216933333333333; H=0.66508; I=0.08800838744757017; J=0.08800838744757018; K=3; L=1; M=0.6602006666666667; N=1.980602; O=0.01082096921202439; P=0.01082096921202439; Q=0.1463620864622222; R=1.31725877816; S=0.0009523360509700702; T=0.0009523360509700703},
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$r = 2
foreach ($row in $rowsData) {
    $c = 1
    foreach ($colName in $cols) {
        $ws.Cells.Item($r, $c).Value = $row[$colName]
        $c = $c + 1
    }
    $r = $r + 1
}
